# Add a new "2021" data column (column R) to the report, mirroring the
# formatting of the existing "2020" column (Q), then move the active
# selection the way the author left it (U4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone the formatting of column Q (rows 4-14, the data block) onto the
#    new column R so every new cell picks up the same number format /
#    font / borders as its neighbour in column Q.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Write the 2021 values.
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 99.4
$ws.Range("R6").Value = 98.1
$ws.Range("R7").Value = 99.319469393395053
$ws.Range("R8").Value = 99.442213297634979
$ws.Range("R9").Value = 99.1
$ws.Range("R10").Value = 99.3
$ws.Range("R11").Value = 99.799160124155549
$ws.Range("R12").Value = 99.3
$ws.Range("R13").Value = 99.538370126605429
$ws.Range("R14").Value = 99.765563948945029

# 3) Leave the selection where the author left it.
[void]$ws.Range("U4").Select()
